$d = $word.ActiveDocument
$ps = $d.PageSetup

# Header/footer distance from edge: 283 twips (14.15pt) / 510 twips (25.5pt)
$ps.HeaderDistance = 14.15
$ps.FooterDistance = 25.5

# Turn on a distinct first-page header/footer layout (adds <w:titlePg/>)
$ps.DifferentFirstPageHeaderFooter = $true
